$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value (RowID) from the old garbage numeric value to 202401
$ws.Range("A2").Value = 202401

# Move the active selection from N2 to A2
$ws.Range("A2").Select()
